$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert column B (生年月日 / birthdate) from text strings to real Excel dates
$ws.Range("B2").Value = 34037
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("B3:B48").PasteSpecial(-4122)

$ws.Range("B3").Value = 36292
$ws.Range("B4").Value = 35832
$ws.Range("B5").Value = 34009
$ws.Range("B6").Value = 36383
$ws.Range("B7").Value = 35282
$ws.Range("B8").Value = 35755
$ws.Range("B9").Value = 37046
$ws.Range("B10").Value = 36297
$ws.Range("B11").Value = 34719
$ws.Range("B12").Value = 38234
$ws.Range("B13").Value = 36890
$ws.Range("B14").Value = 37141
$ws.Range("B15").Value = 36776
$ws.Range("B16").Value = 35570
$ws.Range("B17").Value = 34829
$ws.Range("B18").Value = 35650
$ws.Range("B19").Value = 36035
$ws.Range("B20").Value = 35531
$ws.Range("B21").Value = 36705
$ws.Range("B22").Value = 36684
$ws.Range("B23").Value = 36400
$ws.Range("B24").Value = 34648
$ws.Range("B25").Value = 34585
$ws.Range("B26").Value = 35723
$ws.Range("B27").Value = 35486
$ws.Range("B28").Value = 35273
$ws.Range("B29").Value = 36369
$ws.Range("B30").Value = 36852
$ws.Range("B31").Value = 33434
$ws.Range("B32").Value = 35477
$ws.Range("B33").Value = 36735
$ws.Range("B34").Value = 35667
$ws.Range("B35").Value = 31667
$ws.Range("B36").Value = 36048
$ws.Range("B37").Value = 35466
$ws.Range("B38").Value = 36886
$ws.Range("B39").Value = 37303
$ws.Range("B40").Value = 35962
$ws.Range("B41").Value = 34715
$ws.Range("B42").Value = 35457
$ws.Range("B43").Value = 36104
$ws.Range("B44").Value = 37154
$ws.Range("B45").Value = 35719
$ws.Range("B46").Value = 37615
$ws.Range("B47").Value = 37489
$ws.Range("B48").Value = 37189

# Restore the cell selection shown in the target workbook
$ws.Range("B7").Select()
